# ForEachTagTemplate.xlsx update:
#  - bump active tab to the new last sheet
#  - add two new worksheets: "groupBy" and "orderBy" (JETT forEach groupBy/orderBy examples)

$wb = $excel.ActiveWorkbook

# Reference sheet used as a style/format donor ("VertVert", the first sheet,
# whose A1:E3 already carries the title-merge / header / data-row look).
$styleSrc = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the "groupBy" sheet right after "groupCols" (the current last sheet).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$groupBySheet = $wb.Worksheets.Add($null, $lastSheet)
$groupBySheet.Name = "groupBy"

# --- content (write in this exact order so new shared strings land right) ---
$groupBySheet.Range("A1").Value = '<jt:forEach items="${teams}" var="division" groupBy="divisionName">Division: ${division.obj.divisionName}'
$groupBySheet.Range("A3").Value = '<jt:forEach items="${division.items}" var="team">${team.city}'

$groupBySheet.Range("A2").Value = "City"
$groupBySheet.Range("B2").Value = "Name"
$groupBySheet.Range("C2").Value = "Wins"
$groupBySheet.Range("D2").Value = "Losses"
$groupBySheet.Range("E2").Value = "Pct."
$groupBySheet.Range("B3").Value = '${team.name}'
$groupBySheet.Range("C3").Value = '${team.wins}'
$groupBySheet.Range("D3").Value = '${team.losses}'
$groupBySheet.Range("E3").Value = '${team.pct}</jt:forEach></jt:forEach>'

# --- formatting (copy exact look from the VertVert sheet's first block) ---
# Merge BEFORE pasting formats so the pasted borders aren't re-split by Excel's
# automatic "merged range" border handling.
$groupBySheet.Range("A1:E1").Merge()

$styleSrc.Range("A1:E1").Copy()
$groupBySheet.Range("A1:E1").PasteSpecial(-4122)
$styleSrc.Range("A2:E2").Copy()
$groupBySheet.Range("A2:E2").PasteSpecial(-4122)
$styleSrc.Range("A3:E3").Copy()
$groupBySheet.Range("A3:E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$groupBySheet.Columns.Item(1).ColumnWidth = 14.17
$groupBySheet.Columns.Item(2).ColumnWidth = 14.67

$groupBySheet.Range("A1:E3").Select()

# ---------------------------------------------------------------------------
# 2. Add the "orderBy" sheet right after "groupBy".
# ---------------------------------------------------------------------------
$orderBySheet = $wb.Worksheets.Add($null, $groupBySheet)
$orderBySheet.Name = "orderBy"

# --- content (order chosen so new shared strings land at indices 27-30) ---
$orderBySheet.Range("A1").Value = "Division"
$orderBySheet.Range("F2").Value = '${team.pct}</jt:forEach>'
$orderBySheet.Range("A2").Value = '<jt:forEach items="${teams}" var="team" orderBy="divisionName desc;pct">${team.divisionName}'
$orderBySheet.Range("I1").Value = '<jt:forEach items="${teams}" var="division" groupBy="divisionName" orderBy="divisionName desc;pct">Division: ${division.obj.divisionName}'

$orderBySheet.Range("B1").Value = "City"
$orderBySheet.Range("C1").Value = "Name"
$orderBySheet.Range("D1").Value = "Wins"
$orderBySheet.Range("E1").Value = "Losses"
$orderBySheet.Range("F1").Value = "Pct."
$orderBySheet.Range("B2").Value = '${team.city}'
$orderBySheet.Range("C2").Value = '${team.name}'
$orderBySheet.Range("D2").Value = '${team.wins}'
$orderBySheet.Range("E2").Value = '${team.losses}'

$orderBySheet.Range("I3").Value = '<jt:forEach items="${division.items}" var="team">${team.city}'
$orderBySheet.Range("J2").Value = "Name"
$orderBySheet.Range("I2").Value = "City"
$orderBySheet.Range("K2").Value = "Wins"
$orderBySheet.Range("L2").Value = "Losses"
$orderBySheet.Range("M2").Value = "Pct."
$orderBySheet.Range("J3").Value = '${team.name}'
$orderBySheet.Range("K3").Value = '${team.wins}'
$orderBySheet.Range("L3").Value = '${team.losses}'
$orderBySheet.Range("M3").Value = '${team.pct}</jt:forEach></jt:forEach>'

# --- formatting ---
# Left table: plain header row (style like VertVert row2) + data row (row3).
$styleSrc.Range("A2:E2").Copy()
$orderBySheet.Range("A1:E1").PasteSpecial(-4122)
$styleSrc.Range("E2").Copy()
$orderBySheet.Range("F1").PasteSpecial(-4122)

$styleSrc.Range("A3:D3").Copy()
$orderBySheet.Range("A2:D2").PasteSpecial(-4122)
$styleSrc.Range("D3").Copy()
$orderBySheet.Range("E2").PasteSpecial(-4122)
$styleSrc.Range("E3").Copy()
$orderBySheet.Range("F2").PasteSpecial(-4122)

# Right table: mirrors the "groupBy" sheet layout exactly (title-merge, header, data).
$orderBySheet.Range("I1:M1").Merge()

$styleSrc.Range("A1:E1").Copy()
$orderBySheet.Range("I1:M1").PasteSpecial(-4122)
$styleSrc.Range("A2:E2").Copy()
$orderBySheet.Range("I2:M2").PasteSpecial(-4122)
$styleSrc.Range("A3:E3").Copy()
$orderBySheet.Range("I3:M3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$orderBySheet.Columns.Item(1).ColumnWidth = 14.17
$orderBySheet.Columns.Item(2).ColumnWidth = 14.17
$orderBySheet.Columns.Item(3).ColumnWidth = 14.67
$orderBySheet.Columns.Item(9).ColumnWidth = 14.17
$orderBySheet.Columns.Item(10).ColumnWidth = 14.67

$orderBySheet.Range("A1").Select()
$orderBySheet.Activate()
